# Apply documentation fixes to the "Inventory" worksheet of the
# lux-markdown-tables workbook, per commit:
# "Documentation updates: fixed links, typos, versions, etc."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# /base row: "all environments" -> "all tenants"
$ws.Range("F9").Value = "The base configuration directory applicable to all tenants.  It includes the group configuration, main content database, roles, and application servers."

# /base-unsecured row: expanded description
$ws.Range("F11").Value = "Intended for local developer environments.  Defines a local, non-admin user to perform most of deployments with plus some endpoint consumers."

# /ml-config row: replace the "Gradle Properties" link with a
# "Tenant Configuration" link.
$ws.Range("G8").Value = "[Tenant Configuration](/docs/lux-backend-deployment.md#tenant-configuration)"

# /templates row: reworded description, and drop the now-stale
# "JavaScript Template Files" doc link in G22.
$ws.Range("F22").Value = "Reserved for JavaScript template files used by [/build.gradle](/build.gradle)."
$ws.Range("G22").ClearContents()
